# atualização planilha - 02/08/2024
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dados")
$ws.Activate()

# Copy formatting from row 24 down into the new row 26 (the appended
# entry picks up the same style pattern as the "CAMPO" row above it).
$ws.Range("A24:K24").Copy()
$ws.Range("A26:K26").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# New expense record: CAMPO / 02-08-2024 / "-" / RONALDO / CUSTO / MÃO DE OBRA / und / 1 x 7000
$ws.Range("A26").Value = "CAMPO"
$ws.Range("B26").Value = 45506
$ws.Range("C26").Value = "-"
$ws.Range("D26").Value = "RONALDO"
$ws.Range("E26").Value = "CUSTO"
$ws.Range("F26").Value = "MÃO DE OBRA"
$ws.Range("G26").Value = "und"
$ws.Range("H26").Value = 1
$ws.Range("I26").Value = 7000
$ws.Range("J26").Formula = "=H26*I26"
$ws.Range("K26").Value = "-"

# Column width adjustments (D widened to fit the longer "RONALDO"
# supplier name, K widened slightly)
$ws.Columns.Item(4).ColumnWidth = 31.16
$ws.Columns.Item(11).ColumnWidth = 14.3

# View state: scroll down so row 7 is at the top, select F13
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F13").Select()
